# "week01 of IT security presentation finished"
# Restructure Munka1: drop the old J:K helper list, fix the Lesson01 label,
# and lay out the two new "lesson plan" boxes (rows 6-17 and 18-29).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Stamp the borders/fill/font of the new boxed tables BEFORE we touch
#    the old J:K helper list - we reuse its three row-styles (header /
#    body / bottom) as templates so the new boxes share the exact same
#    formatting, then remove the old list once it has been "harvested".
# ------------------------------------------------------------------
$ws.Range("J2").Copy()
$ws.Range("B6:H6").PasteSpecial(-4122)
$ws.Range("B18:G18").PasteSpecial(-4122)

$ws.Range("J3").Copy()
$ws.Range("B7:H16").PasteSpecial(-4122)
$ws.Range("B19:G28").PasteSpecial(-4122)

$ws.Range("J13").Copy()
$ws.Range("B17:H17").PasteSpecial(-4122)
$ws.Range("B29:G29").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Now the old helper table (content + formatting) can go.
$ws.Range("J2:K13").Clear()

# ------------------------------------------------------------------
# 2) Fix the Lesson01 label in row 3's former spot - it now lives in the
#    new box header, with corrected text ("Lesson0102" -> "Lesson01").
# ------------------------------------------------------------------
$ws.Range("B6").Value = "Lesson01 - Mi az az internet?"
$ws.Range("C6").Value = "Lesson03 - Mi az a wifi?"
$ws.Range("D6").Value = "Lesson05 - Email, internetes levelezés"
$ws.Range("E6").Value = "Lesson07 - Free wifi veszélyei"
$ws.Range("F6").Value = "Lesson09 - Jelszavak működése, titkosítás alapok"
$ws.Range("G6").Value = "Lesson11 - Kvantumszámítógépek, a jövő gépei, titkosításai"
$ws.Range("H6").Value = "Lesson13"

# Body rows 7-16 (left column B = new sub-topic text, column C = old
# Week01 detail list, D-H blank placeholders for future weeks)
$ws.Range("B7").Value = "Kezdete - arpanet"
$ws.Range("B8").Value = "Egyetemi hálózat - > világháló"
$ws.Range("B9").Value = "Lényege régen vs. Ma"
$ws.Range("B10").Value = "Földi összeköttetések - Atlanti kábel, ausztrália kábel"
$ws.Range("B11").Value = "Műholdas összeköttetések - Starlink"
$ws.Range("B12").Value = "WWW - World Wide Web"

$ws.Range("C7").Value = "UTP kábel mint adatátviteli közeg"
$ws.Range("C8").Value = "megmutatni mit tud egy utp kábel"
$ws.Range("C9").Value = "wifi mint adatáviteli közeg"
$ws.Range("C10").Value = "hogyan működik a wifi"
$ws.Range("C11").Value = "rádióhullámok"
$ws.Range("C12").Value = "különböző frekvencia(fontos hogy ne legyen interferencia, pl emergency services)"
$ws.Range("C13").Value = "2.4ghz és 5ghz - unlicensed"
$ws.Range("C14").Value = "2.4ghz és 5ghz közötti különbség"
$ws.Range("C15").Value = "wifi 6"
$ws.Range("C16").Value = "wifi 7"
$ws.Range("C17").Value = "bluetooth vs wifi"

# ------------------------------------------------------------------
# 3) Second lesson-plan box: header row 18, body rows 19-28
# ------------------------------------------------------------------
$ws.Range("B18").Value = "Lesson02 - Mi az az internet 2?"
$ws.Range("C18").Value = "Lesson04 - Áttekintés"
$ws.Range("D18").Value = "Lesson06 - Email átverések, spam, phising, blackmail, data theft"
$ws.Range("E18").Value = "Lesson08 - Áttekintés"
$ws.Range("F18").Value = "Lesson10 - Jelszókezelés, jelszavak feltörése/megszerzése"
$ws.Range("G18").Value = "Lesson12 - Áttekintés"

$ws.Range("B19").Value = "Computer network types - NFC, PAN, LAN, WLAN, WAN, MAN, INTERNET"
$ws.Range("B20").Value = "IP címek - ipv4 minta, ipv6 minta"
$ws.Range("B21").Value = "Mi az a MAC cím = DNS az embernél."
$ws.Range("B22").Value = "Internetet használó eszközök ma"
$ws.Range("B23").Value = "Hány darab eszköz van ami tud csatlakozni az internetre kb"

# ------------------------------------------------------------------
# 4) Column widths for the new, much wider text columns
# ------------------------------------------------------------------
$ws.Columns("B:D").ColumnWidth = 75.85546875
$ws.Columns("E:E").ColumnWidth = 58.42578125
$ws.Columns("F:G").ColumnWidth = 53.85546875
$ws.Columns("H:H").ColumnWidth = 54.7109375

# ------------------------------------------------------------------
# 5) Reset view back to the top-left and select B27, like the author
#    left it after finishing this pass.
# ------------------------------------------------------------------
$ws.Range("B27").Select()
